$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.816.43"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "2.901.50"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'569.37"
$ws.Range("E5").Value = "  -4.41%  "
$ws.Range("D6").Value = "'143.05"
$ws.Range("E6").Value = "  -3.35%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "2.899.10"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("E10").Value = "  -9.15%  "
$ws.Range("E11").Value = "  -3.02%  "
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("D13").Value = "'0.0000233"
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("D14").Value = "'32.14"
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "3.379.21"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").Value = "61.724.80"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "'6.66"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "2.906.07"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "'435.37"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").Value = "'13.29"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "'0.657"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").Value = "'6.92"
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("D24").Value = "'79.70"
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("D25").Value = "'11.83"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").Value = "'10.21"
$ws.Range("E26").Value = "  -9.28%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -5.68%  "
$ws.Range("E29").Value = "  +6.44%  "
$ws.Range("E30").Value = "  -3.77%  "
$ws.Range("E31").Value = "  -4.90%  "
$ws.Range("D32").Value = "'2.07"
$ws.Range("E32").Value = "  -5.05%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.107"
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "'25.66"
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("D36").Value = "'0.962"
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("D37").Value = "'5.45"
$ws.Range("E37").Value = "  -3.93%  "
$ws.Range("D38").Value = "'49.03"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("E39").Value = "  -6.38%  "
$ws.Range("D40").Value = "'2.84"
$ws.Range("E40").Value = "  -10.30%  "
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("D42").Value = "'8.29"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").Value = "'39.45"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").Value = "'0.270"
$ws.Range("E44").Value = "  -5.48%  "
$ws.Range("D45").Value = "2.700.53"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").Value = "'133.67"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").Value = "'0.0335"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'338.60"
$ws.Range("E49").Value = "  -7.16%  "
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").Value = "'21.74"
$ws.Range("E51").Value = "  -5.73%  "
